$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.319.73'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.419.00'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.16%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '568.31'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '179.76'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.07%  '
$ws.Range('E7').Value = '  +1.50%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.410.69'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.24%  '
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.177'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.92%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.640'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.82%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.76'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.22%  '
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('E14').Value = '  +2.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.965.04'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.96%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.33'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.60%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.120'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.01%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.412.03'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '66.202.41'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.00'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.81%  '
$ws.Range('E21').Value = '  +1.90%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '465.25'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.97%  '
$ws.Range('E23').Value = '  +0.63%  '
$ws.Range('E24').Value = '  +8.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.14'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '89.96'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.61%  '
$ws.Range('E27').Value = '  +1.73%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.82'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.89'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.38'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.98%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.93'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.59'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '585.91'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.77%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '62.65'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.13%  '
$ws.Range('E35').Value = '  +1.47%  '
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('E37').Value = '  +4.19%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.59'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.91%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.46'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.30%  '
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.383'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.80%  '
$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0₃0765'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.123.52'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.94'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.68%  '
$ws.Range('E44').Value = '  +2.08%  '
$ws.Range('E45').Value = '  +2.44%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.21'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.72%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.135'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.65'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +14.37%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.999'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '140.76'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.04%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.59'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.13%  '
